$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.789.67"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "3.316.54"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'255.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").Value = "'629.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("D7").Value = "'1.45"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +20.11%  "
$ws.Range("D8").Value = "'0.409"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.95%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "'0.999"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +22.84%  "
$ws.Range("D11").Value = "3.314.62"
$ws.Range("E11").Value = "  -1.28%  "
$ws.Range("E12").Value = "  +3.09%  "
$ws.Range("D13").Value = "'43.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +20.93%  "
$ws.Range("D14").Value = "98.515.92"
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").Value = "'0.0000251"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.60%  "
$ws.Range("D16").Value = "3.939.66"
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("D17").Value = "'5.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("D18").Value = "3.315.46"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").Value = "'16.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.70%  "
$ws.Range("E20").Value = "  -4.75%  "
$ws.Range("D21").Value = "'6.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.55%  "
$ws.Range("D22").Value = "'486.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "'9.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.93%  "
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("D25").Value = "'6.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.23%  "
$ws.Range("D26").Value = "'0.341"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +34.03%  "
$ws.Range("D27").Value = "'90.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.37%  "
$ws.Range("D28").Value = "'12.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.26%  "
$ws.Range("D29").Value = "3.493.62"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").Value = "'0.148"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +16.68%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'11.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +20.06%  "
$ws.Range("B33").Value = "Cronos"
$ws.Range("C33").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D33").Value = "'0.191"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.64%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").Value = "'28.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.78%  "
$ws.Range("D36").Value = "'0.486"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.87%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D37").Value = "'7.45"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.20%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.151"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("D39").Value = "'1.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.04%  "
$ws.Range("D40").Value = "'502.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.97%  "
$ws.Range("D41").Value = "'24.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").Value = "'3.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("D43").Value = "'1.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").Value = "'0.801"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.41%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("E46").Value = "  -2.89%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'1.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.98%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'160.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").Value = "'7.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +16.00%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D50").Value = "'4.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.24%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.857"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.66%  "
